$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Hydrogen - clear the Non-metallic minerals (D) value
$ws.Range("D3").Value = $null

# Row 4: Methanol - correct the Chemicals (C) value
$ws.Range("C4").Value = 0

# Row 5: Ammonia - correct the Chemicals (C) value
$ws.Range("C5").Value = 3612.120285859941

# Row 7: rename "Other" -> "Biogas" and correct its (D) value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 288.3335907277477

# Add new row 8 "Other", cloning row 7's formatting so the new label
# keeps the same bold/centered/boxed style used by the other row headers.
$ws.Range("A7:D7").Copy($ws.Range("A8:D8"))
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 1084.481153000673
